$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 435.26666
$ws.Range("I8").Value = 158
$ws.Range("J8").Value = 504.58334
$ws.Range("K8").Value = 474
$ws.Range("L8").Value = 1513.75002
$ws.Range("M8").Value = -335
$ws.Range("N8").Value = -1791.75002
$ws.Range("H98").Value = 28331.379
$ws.Range("I98").Value = 1015.8571
$ws.Range("J98").Value = 100034.625
$ws.Range("K98").Value = 1015.8571
$ws.Range("L98").Value = 100034.625
$ws.Range("M98").Value = 482.1429000000001
$ws.Range("N98").Value = -103030.625
$ws.Range("H113").Value = 1939.2222
$ws.Range("I113").Value = 1950
$ws.Range("J113").Value = 1917.6666
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 1917.6666
$ws.Range("M113").Value = 1304
$ws.Range("N113").Value = -8425.6666
$ws.Range("H122").Value = 28331.379
$ws.Range("I122").Value = 1015.8571
$ws.Range("J122").Value = 100034.625
$ws.Range("K122").Value = 3047.5713
$ws.Range("L122").Value = 300103.875
$ws.Range("M122").Value = -597.5712999999996
$ws.Range("N122").Value = -305003.875
$ws.Range("H131").Value = 3575
$ws.Range("I131").Value = 3990
$ws.Range("J131").Value = 3461.818
$ws.Range("K131").Value = 11970
$ws.Range("L131").Value = 10385.454
$ws.Range("M131").Value = -6930
$ws.Range("N131").Value = -20465.454
$ws.Range("H137").Value = 2963767.5
$ws.Range("I137").Value = 6993996
$ws.Range("K137").Value = 20981988
$ws.Range("M137").Value = -20979438
$ws.Range("H138").Value = 1471.8586
$ws.Range("I138").Value = 784.8125
$ws.Range("J138").Value = 1800
$ws.Range("K138").Value = 2354.4375
$ws.Range("L138").Value = 5400
$ws.Range("M138").Value = 2785.5625
$ws.Range("N138").Value = -15680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H32").Value = 9889.922
$ws.Range("I32").Value = 8879.421
$ws.Range("J32").Value = 18118.285
$ws.Range("K32").Value = 8879.421
$ws.Range("L32").Value = 18118.285
$ws.Range("M32").Value = -8592.421
$ws.Range("N32").Value = -18692.285
$ws.Range("H61").Value = 2664.4194
$ws.Range("I61").Value = 1927
$ws.Range("J61").Value = 3451
$ws.Range("K61").Value = 1927
$ws.Range("L61").Value = 3451
$ws.Range("M61").Value = -1715
$ws.Range("N61").Value = -3875
$ws.Range("H74").Value = 1877.5349
$ws.Range("I74").Value = 1709
$ws.Range("J74").Value = 2226.6428
$ws.Range("K74").Value = 1709
$ws.Range("L74").Value = 2226.6428
$ws.Range("M74").Value = -835
$ws.Range("N74").Value = -3974.6428
$ws.Range("H77").Value = 1877.5349
$ws.Range("I77").Value = 1709
$ws.Range("J77").Value = 2226.6428
$ws.Range("K77").Value = 8545
$ws.Range("L77").Value = 11133.214
$ws.Range("M77").Value = -4177
$ws.Range("N77").Value = -19869.214
$ws.Range("H136").Value = 2664.4194
$ws.Range("I136").Value = 1927
$ws.Range("J136").Value = 3451
$ws.Range("K136").Value = 5781
$ws.Range("L136").Value = 10353
$ws.Range("M136").Value = -3231
$ws.Range("N136").Value = -15453

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 23500
$ws.Range("I30").Value = 2000
$ws.Range("K30").Value = 2000
$ws.Range("M30").Value = -1875
$ws.Range("H134").Value = 2436.2354
$ws.Range("I134").Value = 1825.125
$ws.Range("J134").Value = 3902.9
$ws.Range("K134").Value = 5475.375
$ws.Range("L134").Value = 11708.7
$ws.Range("M134").Value = -2940.375
$ws.Range("N134").Value = -16778.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6541290.5
$ws.Range("I31").Value = 1900.8667
$ws.Range("J31").Value = 15883276
$ws.Range("K31").Value = 1900.8667
$ws.Range("L31").Value = 15883276
$ws.Range("M31").Value = -1605.8667
$ws.Range("N31").Value = -15883866
$ws.Range("H34").Value = 6541290.5
$ws.Range("I34").Value = 1900.8667
$ws.Range("J34").Value = 15883276
$ws.Range("K34").Value = 1900.8667
$ws.Range("L34").Value = 15883276
$ws.Range("M34").Value = -1698.8667
$ws.Range("N34").Value = -15883680
$ws.Range("H58").Value = 1890.7838
$ws.Range("I58").Value = 1192.8695
$ws.Range("J58").Value = 3037.3572
$ws.Range("K58").Value = 1192.8695
$ws.Range("L58").Value = 3037.3572
$ws.Range("M58").Value = -989.8695
$ws.Range("N58").Value = -3443.3572
$ws.Range("H132").Value = 484324.06
$ws.Range("I132").Value = 986.5
$ws.Range("J132").Value = 2003385
$ws.Range("K132").Value = 2959.5
$ws.Range("L132").Value = 6010155
$ws.Range("M132").Value = -429.5
$ws.Range("N132").Value = -6015215
$ws.Range("H134").Value = 647504.0600000001
$ws.Range("I134").Value = 784065.2
$ws.Range("J134").Value = 237820.67
$ws.Range("K134").Value = 2352195.6
$ws.Range("L134").Value = 713462.01
$ws.Range("M134").Value = -2349660.6
$ws.Range("N134").Value = -718532.01
$ws.Range("H136").Value = 1890.7838
$ws.Range("I136").Value = 1192.8695
$ws.Range("J136").Value = 3037.3572
$ws.Range("K136").Value = 3578.6085
$ws.Range("L136").Value = 9112.071599999999
$ws.Range("M136").Value = -1028.6085
$ws.Range("N136").Value = -14212.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4594.96
$ws.Range("I5").Value = 5483.75
$ws.Range("J5").Value = 1039.8
$ws.Range("K5").Value = 16451.25
$ws.Range("L5").Value = 3119.4
$ws.Range("M5").Value = -16339.25
$ws.Range("N5").Value = -3343.4
$ws.Range("H26").Value = 480
$ws.Range("I26").Value = 475
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 1425
$ws.Range("L26").Value = 1500
$ws.Range("M26").Value = -1137
$ws.Range("N26").Value = -2076
$ws.Range("H122").Value = 4677.769
$ws.Range("I122").Value = 462.6154
$ws.Range("K122").Value = 4163.5386
$ws.Range("M122").Value = -1713.5386
$ws.Range("H135").Value = 4594.96
$ws.Range("I135").Value = 5483.75
$ws.Range("J135").Value = 1039.8
$ws.Range("K135").Value = 49353.75
$ws.Range("L135").Value = 9358.199999999999
$ws.Range("M135").Value = -46818.75
$ws.Range("N135").Value = -14428.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 9750
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 9750
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 9750
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -10096
$ws.Range("H30").Value = 9750
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 9750
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 9750
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -9960
$ws.Range("H107").Value = 7595.2
$ws.Range("I107").Value = 1200
$ws.Range("J107").Value = 9194
$ws.Range("K107").Value = 1200
$ws.Range("L107").Value = 9194
$ws.Range("M107").Value = 720
$ws.Range("N107").Value = -13034
$ws.Range("H126").Value = 6087.48
$ws.Range("I126").Value = 11102.091
$ws.Range("J126").Value = 2147.4285
$ws.Range("K126").Value = 33306.273
$ws.Range("L126").Value = 6442.2855
$ws.Range("M126").Value = -30836.273
$ws.Range("N126").Value = -11382.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4808.763
$ws.Range("I46").Value = 998.1539
$ws.Range("J46").Value = 6790.28
$ws.Range("K46").Value = 998.1539
$ws.Range("L46").Value = 6790.28
$ws.Range("M46").Value = -810.1539
$ws.Range("N46").Value = -7166.28
$ws.Range("H55").Value = 718.6818
$ws.Range("I55").Value = 693.4666999999999
$ws.Range("J55").Value = 772.7143
$ws.Range("K55").Value = 693.4666999999999
$ws.Range("L55").Value = 772.7143
$ws.Range("M55").Value = -520.4666999999999
$ws.Range("N55").Value = -1118.7143
$ws.Range("H68").Value = 3433.6667
$ws.Range("I68").Value = 2901
$ws.Range("J68").Value = 3700
$ws.Range("K68").Value = 2901
$ws.Range("L68").Value = 3700
$ws.Range("M68").Value = -2152
$ws.Range("N68").Value = -5198
$ws.Range("H71").Value = 3433.6667
$ws.Range("I71").Value = 2901
$ws.Range("J71").Value = 3700
$ws.Range("K71").Value = 14505
$ws.Range("L71").Value = 18500
$ws.Range("M71").Value = -10761
$ws.Range("N71").Value = -25988
$ws.Range("H132").Value = 4814
$ws.Range("I132").Value = 3250.5
$ws.Range("J132").Value = 5986.625
$ws.Range("K132").Value = 9751.5
$ws.Range("L132").Value = 17959.875
$ws.Range("M132").Value = -7221.5
$ws.Range("N132").Value = -23019.875
$ws.Range("H136").Value = 1966.8518
$ws.Range("I136").Value = 1586.8636
$ws.Range("K136").Value = 4760.5908
$ws.Range("M136").Value = -2210.5908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1360617.6
$ws.Range("I132").Value = 1740427
$ws.Range("J132").Value = 4155.857
$ws.Range("K132").Value = 5221281
$ws.Range("L132").Value = 12467.571
$ws.Range("M132").Value = -5218751
$ws.Range("N132").Value = -17527.571
$ws.Range("H136").Value = 477407.75
$ws.Range("I136").Value = 834345.6
$ws.Range("J136").Value = 1490.5834
$ws.Range("K136").Value = 2503036.8
$ws.Range("L136").Value = 4471.7502
$ws.Range("M136").Value = -2500486.8
$ws.Range("N136").Value = -9571.7502
